$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 619.25
$ws.Range("I18").Value = 619.25
$ws.Range("K18").Value = 619.25
$ws.Range("M18").Value = -335.25

$ws.Range("H28").Value = 4231.6665
$ws.Range("I28").Value = 6917.2
$ws.Range("J28").Value = 874.75
$ws.Range("K28").Value = 6917.2
$ws.Range("L28").Value = 874.75
$ws.Range("M28").Value = -6432.2
$ws.Range("N28").Value = -1844.75

$ws.Range("H43").Value = 12920.333
$ws.Range("I43").Value = 20962.5
$ws.Range("K43").Value = 20962.5
$ws.Range("M43").Value = -20893.5

$ws.Range("H62").Value = 26673448
$ws.Range("I62").Value = 26673448
$ws.Range("K62").Value = 26673448
$ws.Range("M62").Value = -26672824

$ws.Range("H65").Value = 26673448
$ws.Range("I65").Value = 26673448
$ws.Range("K65").Value = 133367240
$ws.Range("M65").Value = -133364120

$ws.Range("H76").Value = 6255699.5
$ws.Range("I76").Value = 8338778.5
$ws.Range("J76").Value = 6462
$ws.Range("K76").Value = 8338778.5
$ws.Range("L76").Value = 6462
$ws.Range("M76").Value = -8338463.5
$ws.Range("N76").Value = -7092

$ws.Range("H79").Value = 6255699.5
$ws.Range("I79").Value = 8338778.5
$ws.Range("J79").Value = 6462
$ws.Range("K79").Value = 8338778.5
$ws.Range("L79").Value = 6462
$ws.Range("M79").Value = -8337686.5
$ws.Range("N79").Value = -8646

$ws.Range("H92").Value = 1071.1578
$ws.Range("I92").Value = 1266.0625
$ws.Range("K92").Value = 1266.0625
$ws.Range("M92").Value = -18.0625

$ws.Range("H107").Value = 62501244
$ws.Range("I107").Value = 62501244
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 62501244
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws.Range("H137").Value = 27787458
$ws.Range("I137").Value = 62501076
$ws.Range("K137").Value = 187503228
$ws.Range("M137").Value = -187500678

$ws.Range("H138").Value = 2656.0444
$ws.Range("I138").Value = 1862.3077
$ws.Range("K138").Value = 5586.9231
$ws.Range("M138").Value = -446.9231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25447.682
$ws.Range("I32").Value = 27060.045
$ws.Range("K32").Value = 27060.045
$ws.Range("M32").Value = -26773.045

$ws.Range("H45").Value = 17335.334
$ws.Range("I45").Value = 24503
$ws.Range("K45").Value = 24503
$ws.Range("M45").Value = -24126

$ws.Range("H61").Value = 8277.25
$ws.Range("I61").Value = 7072
$ws.Range("K61").Value = 7072
$ws.Range("M61").Value = -6860

$ws.Range("H97").Value = 3368114
$ws.Range("I97").Value = 4630832
$ws.Range("K97").Value = 4630832
$ws.Range("M97").Value = -4630336

$ws.Range("H132").Value = 3682.3
$ws.Range("I132").Value = 3376.75
$ws.Range("J132").Value = 4904.5
$ws.Range("K132").Value = 10130.25
$ws.Range("L132").Value = 14713.5
$ws.Range("M132").Value = -7600.25
$ws.Range("N132").Value = -19773.5

$ws.Range("H136").Value = 8277.25
$ws.Range("I136").Value = 7072
$ws.Range("K136").Value = 21216
$ws.Range("M136").Value = -18666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 644271.3
$ws.Range("I86").Value = 2474.75
$ws.Range("K86").Value = 2474.75
$ws.Range("M86").Value = -1351.75

$ws.Range("H89").Value = 644271.3
$ws.Range("I89").Value = 2474.75
$ws.Range("K89").Value = 12373.75
$ws.Range("M89").Value = -6757.75

$ws.Range("H134").Value = 13789.066
$ws.Range("I134").Value = 13666.909
$ws.Range("K134").Value = 41000.727
$ws.Range("M134").Value = -38465.727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 9670.532999999999
$ws.Range("I134").Value = 9654.923000000001
$ws.Range("K134").Value = 28964.769
$ws.Range("M134").Value = -26429.769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 28185380
$ws.Range("I4").Value = 42268104
$ws.Range("K4").Value = 126804312
$ws.Range("M4").Value = -126804200

$ws.Range("H12").Value = 3235.889
$ws.Range("J12").Value = 3015.625
$ws.Range("L12").Value = 9046.875
$ws.Range("N12").Value = -9392.875

$ws.Range("H103").Value = 5259.625
$ws.Range("I103").Value = 5349.6665
$ws.Range("J103").Value = 5205.6
$ws.Range("K103").Value = 16048.9995
$ws.Range("L103").Value = 15616.8
$ws.Range("M103").Value = -15169.9995
$ws.Range("N103").Value = -17374.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 11545.272
$ws.Range("I33").Value = 8333.333000000001
$ws.Range("J33").Value = 25999
$ws.Range("K33").Value = 8333.333000000001
$ws.Range("L33").Value = 25999
$ws.Range("M33").Value = -8081.333000000001
$ws.Range("N33").Value = -26503

$ws.Range("H132").Value = 8877.625
$ws.Range("I132").Value = 8003
$ws.Range("K132").Value = 24009
$ws.Range("M132").Value = -21479

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3300.7632
$ws.Range("I22").Value = 2027.8182
$ws.Range("J22").Value = 5051.0625
$ws.Range("K22").Value = 2027.8182
$ws.Range("L22").Value = 5051.0625
$ws.Range("M22").Value = -1732.8182
$ws.Range("N22").Value = -5641.0625

$ws.Range("H27").Value = 3300.7632
$ws.Range("I27").Value = 2027.8182
$ws.Range("J27").Value = 5051.0625
$ws.Range("K27").Value = 2027.8182
$ws.Range("L27").Value = 5051.0625
$ws.Range("M27").Value = -1920.8182
$ws.Range("N27").Value = -5265.0625

$ws.Range("H132").Value = 14936.875
$ws.Range("I132").Value = 14999.333
$ws.Range("K132").Value = 44997.999
$ws.Range("M132").Value = -42467.999

$ws.Range("H136").Value = 4078.5
$ws.Range("I136").Value = 2081.4736
$ws.Range("K136").Value = 6244.4208
$ws.Range("M136").Value = -3694.4208

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 28028.5
$ws.Range("I61").Value = 30000
$ws.Range("K61").Value = 30000
$ws.Range("M61").Value = -29708

$ws.Range("H62").Value = 16756.916
$ws.Range("I62").Value = 20688.555
$ws.Range("J62").Value = 4962
$ws.Range("K62").Value = 20688.555
$ws.Range("L62").Value = 4962
$ws.Range("M62").Value = -20064.555
$ws.Range("N62").Value = -6210

$ws.Range("H65").Value = 16756.916
$ws.Range("I65").Value = 20688.555
$ws.Range("J65").Value = 4962
$ws.Range("K65").Value = 103442.775
$ws.Range("L65").Value = 24810
$ws.Range("M65").Value = -100322.775
$ws.Range("N65").Value = -31050

$ws.Range("H136").Value = 4114.5938
$ws.Range("I136").Value = 2856.5833
$ws.Range("K136").Value = 8569.749899999999
$ws.Range("M136").Value = -6019.749899999999
